$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row -> new D (Price) and new E (Volume(1h)) values.
# Only rows/columns actually changed per the diff are included (D is
# omitted for rows where the price did not change).
# Values are plain text in the sheet (e.g. "261.03", "1.88%"), so force
# the cell format to Text before assigning to avoid Excel auto-converting
# them to numbers/percentages.
$updates = @(
    @{ Row = 2;  D = "261.03";     E = "1.88%" },
    @{ Row = 3;  D = "27.23";      E = "1.47%" },
    @{ Row = 4;  D = "4.746";      E = "5.10%" },
    @{ Row = 5;  D = "0.06085";    E = "3.48%" },
    @{ Row = 6;  D = "6.667";      E = "0.88%" },
    @{ Row = 7;  D = "0.8461";     E = "-0.49%" },
    @{ Row = 8;  D = "0.9229";     E = "-0.45%" },
    @{ Row = 9;  D = "0.1408";     E = "2.31%" },
    @{ Row = 10; D = "0.04966";    E = "10.11%" },
    @{ Row = 11; D = "0.07103";    E = "0.71%" },
    @{ Row = 12; D = "0.03132";    E = "2.21%" },
    @{ Row = 13; D = "0.09080";    E = "-0.06%" },
    @{ Row = 14; D = "0.001529";   E = "-0.43%" },
    @{ Row = 15; D = "0.0006084";  E = "0.25%" },
    @{ Row = 16; D = "0.006140";   E = "0.54%" },
    @{ Row = 17; D = "3.452";      E = "-0.75%" },
    @{ Row = 18; D = "3.149";      E = "-0.66%" },
    @{ Row = 20; E = "2.59%" },
    @{ Row = 22; D = "4.093";      E = "4.67%" },
    @{ Row = 23; D = "0.04235";    E = "-0.56%" },
    @{ Row = 24; D = "0.001219";   E = "-0.19%" },
    @{ Row = 26; E = "-0.01%" },
    @{ Row = 27; D = "0.0001575";  E = "3.40%" },
    @{ Row = 40; D = "0.03872";    E = "1.92%" },
    @{ Row = 41; E = "1.43%" },
    @{ Row = 42; E = "-33.92%" },
    @{ Row = 44; E = "-9.51%" },
    @{ Row = 45; D = "0.00005319"; E = "-0.41%" },
    @{ Row = 46; E = "0.03%" },
    @{ Row = 47; E = "1.28%" },
    @{ Row = 48; E = "-46.44%" },
    @{ Row = 49; E = "0.03%" },
    @{ Row = 50; E = "0.03%" }
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($u.ContainsKey("D")) {
        $cell = $ws.Cells.Item($r, 4)
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
    }
    if ($u.ContainsKey("E")) {
        $cell = $ws.Cells.Item($r, 5)
        $cell.NumberFormat = "@"
        $cell.Value = $u.E
    }
}
